$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" text in cell A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$oldText = $wsHoja1.Range("A1").Value()
$newText = $oldText -replace [regex]::Escape("1000 Bs = 1.82 = 6743.17 pesos"), "1000 Bs = 1.87 = 6911.32 pesos"
$newText = $newText -replace [regex]::Escape("6743.17 pesos = 1.81 = 822.79 Bs"), "6911.32 pesos = 1.85 = 924.96 Bs"
$wsHoja1.Range("A1").Value = $newText

# --- tasas: update rate cells ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 535.996
$wsTasas.Range("O10").Value = 3704.44
$wsTasas.Range("N12").Value = 3736
$wsTasas.Range("O12").Value = 500
